# commit: add calendar bug ver & fix triangle bug
#
# - Clear the stray "actual value" / "execution result" columns (F:G) that had
#   been populated for every data row (rows 2-15) on Sheet1; only the header
#   row (row 1) keeps its F/G labels.
# - Fix the "triangle bug": row 15's expected-value cell (E15) held the
#   string "0" instead of the correctly formatted "0.0".
# - Leave the cursor/selection on I9 (matches the saved view state).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the leftover "actual value"/"result" data that shouldn't have been
# written for the data rows (header row F1/G1 stays intact).
$ws.Range("F2:G15").ClearContents()

# Correct the expected value for the last test case row.
$ws.Range("E15").Value = "0.0"

# Restore the saved selection/active cell.
$ws.Range("I9").Select()
